$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows before the current row 294, shifting existing rows 294-365 down to 299-370
$ws.Rows("294:298").Insert()

# Populate the 5 newly inserted rows (294-298) with the new data block
$ws.Range("A294").Value = 6
$ws.Range("B294").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C294").Value = "Metropolitana"
$ws.Range("D294").Value = 44463
$ws.Range("E294").Value = 13
$ws.Range("F294").Value = 100112013
$ws.Range("G294").Value = "Alcachofa"
$ws.Range("H294").Value = "Argentina(o)"
$ws.Range("I294").Value = "Primera"
$ws.Range("J294").Value = 400
$ws.Range("K294").Value = 8000
$ws.Range("L294").Value = 10000
$ws.Range("M294").Value = 9150
$ws.Range("N294").Value = "$/caja 50 unidades"
$ws.Range("O294").Value = "Provincia de Limarí"
$ws.Range("P294").Value = 183
$ws.Range("Q294").Value = 50
$ws.Range("R294").Value = "Hortaliza"
$ws.Range("A295").Value = 6
$ws.Range("B295").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C295").Value = "Metropolitana"
$ws.Range("D295").Value = 44463
$ws.Range("E295").Value = 13
$ws.Range("F295").Value = 100112013
$ws.Range("G295").Value = "Alcachofa"
$ws.Range("H295").Value = "Argentina(o)"
$ws.Range("I295").Value = "Segunda"
$ws.Range("J295").Value = 400
$ws.Range("K295").Value = 7000
$ws.Range("L295").Value = 8000
$ws.Range("M295").Value = 7575
$ws.Range("N295").Value = "$/caja 70 unidades"
$ws.Range("O295").Value = "Provincia de Limarí"
$ws.Range("P295").Value = 108
$ws.Range("Q295").Value = 70
$ws.Range("R295").Value = "Hortaliza"
$ws.Range("A296").Value = 6
$ws.Range("B296").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C296").Value = "Metropolitana"
$ws.Range("D296").Value = 44463
$ws.Range("E296").Value = 13
$ws.Range("F296").Value = 100112013
$ws.Range("G296").Value = "Alcachofa"
$ws.Range("H296").Value = "Española"
$ws.Range("I296").Value = "Extra"
$ws.Range("J296").Value = 550
$ws.Range("K296").Value = 10000
$ws.Range("L296").Value = 12000
$ws.Range("M296").Value = 10836
$ws.Range("N296").Value = "$/caja 25 unidades"
$ws.Range("O296").Value = "Provincia de Limarí"
$ws.Range("P296").Value = 10836
$ws.Range("Q296").Value = 1
$ws.Range("R296").Value = "Hortaliza"
$ws.Range("A297").Value = 6
$ws.Range("B297").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C297").Value = "Metropolitana"
$ws.Range("D297").Value = 44463
$ws.Range("E297").Value = 13
$ws.Range("F297").Value = 100112013
$ws.Range("G297").Value = "Alcachofa"
$ws.Range("H297").Value = "Española"
$ws.Range("I297").Value = "Primera"
$ws.Range("J297").Value = 400
$ws.Range("K297").Value = 8000
$ws.Range("L297").Value = 9000
$ws.Range("M297").Value = 8425
$ws.Range("N297").Value = "$/caja 30 unidades"
$ws.Range("O297").Value = "Provincia de Limarí"
$ws.Range("P297").Value = 281
$ws.Range("Q297").Value = 30
$ws.Range("R297").Value = "Hortaliza"
$ws.Range("A298").Value = 6
$ws.Range("B298").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C298").Value = "Metropolitana"
$ws.Range("D298").Value = 44463
$ws.Range("E298").Value = 13
$ws.Range("F298").Value = 100112013
$ws.Range("G298").Value = "Alcachofa"
$ws.Range("H298").Value = "Madrigal"
$ws.Range("I298").Value = "Primera"
$ws.Range("J298").Value = 490
$ws.Range("K298").Value = 7000
$ws.Range("L298").Value = 8000
$ws.Range("M298").Value = 7469
$ws.Range("N298").Value = "$/caja 40 unidades"
$ws.Range("O298").Value = "Región de Coquimbo"
$ws.Range("P298").Value = 187
$ws.Range("Q298").Value = 40
$ws.Range("R298").Value = "Hortaliza"
